{"js": "// Append a continuation run to the \"Result: 3223 images\" paragraph, then add\n// three brand-new paragraphs after it, per the diff:\n//   \"Result: 3223 images\" + \": CN group=1782, AD group=1441\"\n//   \"(actually 1780 cn and 1437 ad)\"\n//   \"Number of subjects: 178 AD, 216 CN \u2013 should check after this deleting of bad images?\"\n//   \"2d rbg images: 104961 of CN and 84.783 of AD\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The target paragraph is the last one in the document (\"Result: 3223 images\").\nconst resultParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Extend the existing paragraph's text in place (new trailing run).\nresultParagraph.insertText(\": CN group=1782, AD group=1441\", \"End\");\n\n// Insert the three new paragraphs right after it, in order.\nconst p2 = resultParagraph.insertParagraph(\"(actually 1780 cn and 1437 ad)\", \"After\");\nconst p3 = p2.insertParagraph(\n  \"Number of subjects: 178 AD, 216 CN \u2013 should check after this deleting of bad images?\",\n  \"After\"\n);\np3.insertParagraph(\"2d rbg images: 104961 of CN and 84.783 of AD\", \"After\");\n\nawait context.sync();\n", "ps1": "# Append a continuation run to the \"Result: 3223 images\" paragraph, then add\n# three brand-new paragraphs after it, per the diff:\n#   \"Result: 3223 images\" + \": CN group=1782, AD group=1441\"\n#   \"(actually 1780 cn and 1437 ad)\"\n#   \"Number of subjects: 178 AD, 216 CN \u2013 should check after this deleting of bad images?\"\n#   \"2d rbg images: 104961 of CN and 84.783 of AD\"\n\n$d = $word.ActiveDocument\n\n# The target paragraph is the last one in the document (\"Result: 3223 images\").\n$resultParagraph = $d.Paragraphs.Last\n$resultParagraph.Range.InsertAfter(\": CN group=1782, AD group=1441\")\n\n# Start a new paragraph right after it and fill it in.\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n$d.Paragraphs.Last.Range.InsertAfter(\"(actually 1780 cn and 1437 ad)\")\n\n# Third new paragraph (the en dash is Unicode U+2013).\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n$enDash = [char]0x2013\n$d.Paragraphs.Last.Range.InsertAfter(\"Number of subjects: 178 AD, 216 CN \" + $enDash + \" should check after this deleting of bad images?\")\n\n# Fourth new paragraph.\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n$d.Paragraphs.Last.Range.InsertAfter(\"2d rbg images: 104961 of CN and 84.783 of AD\")\n"}
